# Updates the Coin/Link/Price/Volume(1h) table on Sheet1 to match the
# "Updated symbol list" GitHub Actions commit (17 Dec 2022 18:36 UTC).
#
# Column D ("Price") is stored as text (e.g. "237.40", "0.05580") so that
# trailing/leading zeros survive. A plain $ws.Range(...).Value = "237.40"
# would be auto-coerced to the number 237.4 by Excel, dropping the zero,
# so price cells are written with a leading apostrophe ($apos) to force
# Excel to keep them as literal text, exactly like the source cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$apos = "'"

$ws.Range("D2").Value = $apos + '237.40'

$ws.Range("D3").Value = $apos + '21.93'

$ws.Range("D4").Value = $apos + '5.348'

$ws.Range("D5").Value = $apos + '0.05580'

$ws.Range("D6").Value = $apos + '6.460'

$ws.Range("D7").Value = $apos + '3.355'

$ws.Range("D8").Value = $apos + '0.8003'

$ws.Range("D9").Value = $apos + '1.043'

$ws.Range("D10").Value = $apos + '0.1389'

$ws.Range("D11").Value = $apos + '0.07288'

$ws.Range("D12").Value = $apos + '0.03198'

$ws.Range("D13").Value = $apos + '0.02964'

$ws.Range("D14").Value = $apos + '0.09236'

$ws.Range("D15").Value = $apos + '0.001659'

$ws.Range("D16").Value = $apos + '3.251'

$ws.Range("D17").Value = $apos + '0.04770'

$ws.Range("D18").Value = $apos + '0.0005712'
$ws.Range("E18").Value = '17OneONE'

$ws.Range("D19").Value = $apos + '0.006213'

$ws.Range("D20").Value = $apos + '0.005077'

$ws.Range("D21").Value = $apos + '0.001052'

$ws.Range("D22").Value = $apos + '0.0001499'

$ws.Range("D23").Value = $apos + '0.0003995'

$ws.Range("D24").Value = $apos + '3.940'

$ws.Range("D25").Value = $apos + '2.203'

$ws.Range("D40").Value = $apos + '0.04116'

$ws.Range("D41").Value = $apos + '0.007094'

$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = $apos + '0.1039'
$ws.Range("E42").Value = '41BKEXTokenBKK'

$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = $apos + '0.002939'
$ws.Range("E43").Value = '42CEJICEJI'

$ws.Range("D44").Value = $apos + '0.008929'

$ws.Range("D45").Value = $apos + '0.00005434'

$ws.Range("D47").Value = $apos + '0.6754'

$ws.Range("D48").Value = $apos + '0.03429'
$ws.Range("E48").Value = '47BOLOBOLOWorstin24h'

$ws.Range("D49").Value = $apos + '0.00002101'
